$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New task row (row 13): "Addition of Freshmen" task name cell (B13) plus its
# finish-date cell (D13), matching the pattern used by rows 11/12.
$ws.Range("B12").Copy() | Out-Null
$ws.Range("B13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("D12").Copy() | Out-Null
$ws.Range("D13").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$excel.CutCopyMode = 0

$ws.Range("B13").Value = "Addition of Freshmen"
$ws.Range("D13").Value = "Sun 12/17/17"

# B13 keeps the left/right border but drops the top/bottom rule that B12 has,
# since it is now the last row of the table.
$ws.Range("B13").Borders.Item(8).LineStyle = 0   # xlEdgeTop -> none
$ws.Range("B13").Borders.Item(9).LineStyle = 0   # xlEdgeBottom -> none

$ws.Rows.Item(13).RowHeight = 30

# Refresh view state: scroll so row 5 is the top row and select G10.
$ws.Range("G10").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
